$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The survey response row (row 4, columns B:BZ) was reset from the
# "2" placeholder answers to "0" for every question.
$ws.Range("B4:BZ4").Value = 0

# The author's view/selection moved to the far right of the sheet
# (scrolled so column BO is first visible, with CB9 as the active cell).
$win = $excel.ActiveWindow
try { $win.ScrollColumn = 67 } catch {}
try { $win.ScrollRow = 1 } catch {}
$ws.Range("CB9").Select()
